$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 1207.65
$ws.Range("I41").Value = 1047.5454
$ws.Range("J41").Value = 1403.3334
$ws.Range("K41").Value = 1047.5454
$ws.Range("L41").Value = 1403.3334
$ws.Range("M41").Value = -607.5454
$ws.Range("N41").Value = -2283.3334

# Row 64
$ws.Range("H64").Value = 3622.8386
$ws.Range("J64").Value = 3556.3
$ws.Range("L64").Value = 3556.3
$ws.Range("N64").Value = -4052.3

# Row 67
$ws.Range("H67").Value = 3622.8386
$ws.Range("J67").Value = 3556.3
$ws.Range("L67").Value = 3556.3
$ws.Range("N67").Value = -5272.3

# Row 76
$ws.Range("H76").Value = 3263.4211
$ws.Range("I76").Value = 3189.5557
$ws.Range("J76").Value = 3329.9
$ws.Range("K76").Value = 3189.5557
$ws.Range("L76").Value = 3329.9
$ws.Range("M76").Value = -2874.5557
$ws.Range("N76").Value = -3959.9

# Row 79
$ws.Range("H79").Value = 3263.4211
$ws.Range("I79").Value = 3189.5557
$ws.Range("J79").Value = 3329.9
$ws.Range("K79").Value = 3189.5557
$ws.Range("L79").Value = 3329.9
$ws.Range("M79").Value = -2097.5557
$ws.Range("N79").Value = -5513.9

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1156.05
$ws.Range("I45").Value = 1290.4667
$ws.Range("K45").Value = 1290.4667
$ws.Range("M45").Value = -913.4666999999999

# Row 48
$ws.Range("H48").Value = 183494.67
$ws.Range("J48").Value = 183494.67
$ws.Range("L48").Value = 183494.67
$ws.Range("N48").Value = -184262.67

# Row 102
$ws.Range("H102").Value = 20834082
$ws.Range("I102").Value = 20834082
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 20834082
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -20832460
$ws.Range("N102").ClearContents()

# Row 110
$ws.Range("H110").Value = 2406.4
$ws.Range("I110").Value = 1886.6666
$ws.Range("J110").Value = 2629.1428
$ws.Range("K110").Value = 1886.6666
$ws.Range("L110").Value = 2629.1428
$ws.Range("M110").Value = 158.3334
$ws.Range("N110").Value = -6719.1428

# Row 112
$ws.Range("H112").Value = 15575
$ws.Range("J112").Value = 15575
$ws.Range("L112").Value = 15575
$ws.Range("N112").Value = -18529

$ws = $wb.Worksheets.Item("BSM")
# Row 42
$ws.Range("H42").Value = 132127.2
$ws.Range("J42").Value = 132127.2
$ws.Range("L42").Value = 132127.2
$ws.Range("N42").Value = -132783.2

# Row 43
$ws.Range("H43").Value = 237342
$ws.Range("J43").Value = 237342
$ws.Range("L43").Value = 237342
$ws.Range("N43").Value = -237704

# Row 47
$ws.Range("H47").Value = 149842
$ws.Range("J47").Value = 149842
$ws.Range("L47").Value = 149842
$ws.Range("N47").Value = -150882

# Row 61
$ws.Range("H61").Value = 16000
$ws.Range("J61").Value = 16000
$ws.Range("L61").Value = 16000
$ws.Range("N61").Value = -16626

# Row 100
$ws.Range("H100").Value = 9000
$ws.Range("J100").Value = 9000
$ws.Range("L100").Value = 9000
$ws.Range("N100").Value = -11164

# Row 110
$ws.Range("H110").Value = 18499.75
$ws.Range("J110").Value = 18499.75
$ws.Range("L110").Value = 18499.75
$ws.Range("N110").Value = -26679.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2028.775
$ws.Range("I31").Value = 963.8
$ws.Range("K31").Value = 963.8
$ws.Range("M31").Value = -668.8

# Row 34
$ws.Range("H34").Value = 2028.775
$ws.Range("I34").Value = 963.8
$ws.Range("K34").Value = 963.8
$ws.Range("M34").Value = -761.8

# Row 92
$ws.Range("H92").Value = 44249.75
$ws.Range("J92").Value = 44249.75
$ws.Range("L92").Value = 44249.75
$ws.Range("N92").Value = -49241.75

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1353.5652
$ws.Range("I5").Value = 1353.5652
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4060.6956
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3948.6956
$ws.Range("N5").ClearContents()

# Row 11
$ws.Range("H11").Value = 1966.25
$ws.Range("I11").Value = 2455
$ws.Range("J11").Value = 500
$ws.Range("K11").Value = 7365
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = -7225
$ws.Range("N11").Value = -1780

# Row 34
$ws.Range("H34").Value = 1293.25
$ws.Range("I34").Value = 778.4
$ws.Range("J34").Value = 1527.2727
$ws.Range("K34").Value = 2335.2
$ws.Range("L34").Value = 4581.8181
$ws.Range("M34").Value = -2251.2
$ws.Range("N34").Value = -4749.8181

# Row 39
$ws.Range("H39").Value = 2011.75
$ws.Range("J39").Value = 2011.75
$ws.Range("L39").Value = 6035.25
$ws.Range("N39").Value = -6623.25

# Row 55
$ws.Range("H55").Value = 2993.6365
$ws.Range("J55").Value = 2993.6365
$ws.Range("L55").Value = 8980.9095
$ws.Range("N55").Value = -9334.9095

# Row 122
$ws.Range("H122").Value = 774.3333
$ws.Range("I122").Value = 719.5714
$ws.Range("K122").Value = 6476.1426
$ws.Range("M122").Value = -4026.1426

# Row 124
$ws.Range("H124").Value = 2575
$ws.Range("J124").Value = 3100
$ws.Range("L124").Value = 9300
$ws.Range("N124").Value = -19120

# Row 132
$ws.Range("H132").Value = 957.8
$ws.Range("J132").Value = 995
$ws.Range("L132").Value = 8955
$ws.Range("N132").Value = -14015

# Row 135
$ws.Range("H135").Value = 1353.5652
$ws.Range("I135").Value = 1353.5652
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 12182.0868
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -9647.086800000001
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4551.25
$ws.Range("I80").Value = 905
$ws.Range("K80").Value = 905
$ws.Range("M80").Value = 93

# Row 83
$ws.Range("H83").Value = 4551.25
$ws.Range("I83").Value = 905
$ws.Range("K83").Value = 4525
$ws.Range("M83").Value = 467

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 249.84616
$ws.Range("I55").Value = 99.42856999999999
$ws.Range("K55").Value = 99.42856999999999
$ws.Range("M55").Value = 73.57143000000001

# Row 105
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988

# Row 110
$ws.Range("H110").Value = 17909
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 17909
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 17909
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = -26089

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 26637.5
$ws.Range("J92").Value = 26637.5
$ws.Range("L92").Value = 26637.5
$ws.Range("N92").Value = -31629.5

# Row 136
$ws.Range("H136").Value = 898.75
$ws.Range("I136").Value = 531.6667
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 1595.0001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 954.9999
$ws.Range("N136").Value = -11100
